# Add an "error" (variance) column next to the existing "req" (mean) column,
# and replace several static mean values with AVERAGE(...) formulas computed
# from the underlying replicate readings. Also drop the two now-unused rows
# (A14 "water loss", A15 "drinking water treatment").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a new column I (old I "model" group shifts to J) ---
$ws.Columns.Item(9).Insert()

# New header for the inserted column I ("var" = variance of the req readings)
$ws.Range("I1").Value = "var"

# --- 2. Row 2 (reverse osmosis) ---
$ws.Range("H2").Formula = "=AVERAGE(1, 2.5)"
$ws.Range("I2").Formula = "=VAR.P(1,2.5)"

# --- 3. Row 3 (ultrafiltration) ---
$ws.Range("H3").Formula = "=AVERAGE(0.13, 0.177, 0.201, 0.2, 0.3)"
$ws.Range("I3").Formula = "=VAR.P(0.07, 0.1, 0.2)"

# --- 4. Row 4 (granular activated carbon) -- H4 / J4 stay as-is, no I4 ---

# --- 5. Row 5 (ozonation) -- H5 formula is unchanged, just add the variance ---
$ws.Range("I5").Formula = "=VAR.P(0.05*3.79,0.12*3.79)/24"

# --- 6. Row 6 (uv oxidation) ---
$ws.Range("H6").Formula = "=AVERAGE(36.46, 36.46, 35.61, 36.64, 36.91, 35.43) / 1000"
$ws.Range("I6").Formula = "=VAR.P(36.46, 36.46, 35.61, 36.64, 36.91, 35.43) / 1000"

# --- 7. Row 7 (microfiltration) ---
$ws.Range("H7").Formula = "=AVERAGE(0.11, 0.24, 0.2, 0.2, 0.3)"
$ws.Range("I7").Formula = "=VAR.P(0.11, 0.24, 0.2, 0.2, 0.3)"

# --- 8. Row 8 (brackish water desalination) ---
$ws.Range("H8").Formula = "=AVERAGE(1.02, 2.57)"
$ws.Range("I8").Formula = "=VAR.P(1.02, 2.57)"

# --- 9. Row 9 (seawater desalination) ---
$ws.Range("H9").Formula = "=AVERAGE(2.58,5.47)"
$ws.Range("I9").Formula = "=VAR.P(2.58,5.47)"

# --- 10. Row 10 (groundwater pumping) -- nothing to do ---

# --- 11. Row 11 (coagulation) -- H11 formula unchanged, just add the variance ---
$ws.Range("I11").Formula = "=VAR.P(0.4, 0.7)"

# --- 12. Row 12 (groundwater recharge) -- H12 stays as-is, no I12 ---

# --- 13. Row 13 (nanofiltration) -- brand new req/var pair ---
$ws.Range("H13").Formula = "=AVERAGE(1.33, 0.68, 1.17)"
$ws.Range("I13").Formula = "=VAR.P(1.33, 0.68, 1.17)"

# --- 14. Drop the two trailing rows that no longer belong (A14, A15) ---
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# --- 15. Leave the selection where the author's session ended up ---
$ws.Range("Q25").Select() | Out-Null
